$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cell -> new text value (prefixed with an apostrophe when the text
# looks like a number, e.g. "1.00" or "596.29", so Excel keeps it as
# text instead of silently normalising it to a numeric value).
$updates = [ordered]@{
    "D2" = "67.462.77"
    "E2" = "  -4.17%  "
    "D3" = "3.278.93"
    "E3" = "  -6.22%  "
    "D4" = "'0.999"
    "E4" = "  +0.11%  "
    "D5" = "'596.29"
    "E5" = "  -3.43%  "
    "E6" = "  -10.83%  "
    "D7" = "'1.00"
    "E7" = "  -0.11%  "
    "D8" = "3.271.38"
    "E8" = "  -6.34%  "
    "E9" = "  -9.71%  "
    "D10" = "'0.173"
    "E10" = "  -12.18%  "
    "D11" = "'6.82"
    "E11" = "  -5.11%  "
    "D12" = "'0.510"
    "E12" = "  -11.33%  "
    "D13" = "'38.64"
    "E13" = "  -14.78%  "
    "D14" = "'0.0000247"
    "E14" = "  -9.33%  "
    "D15" = "3.802.47"
    "E15" = "  -6.30%  "
    "D16" = "67.464.03"
    "E16" = "  -4.24%  "
    "D17" = "3.280.99"
    "E17" = "  -5.83%  "
    "D18" = "'536.67"
    "E18" = "  -10.44%  "
    "E19" = "  -5.94%  "
    "E20" = "  -13.27%  "
    "D21" = "'15.17"
    "E21" = "  -13.15%  "
    "E22" = "  -12.29%  "
    "E23" = "  -12.22%  "
    "D24" = "'86.04"
    "E24" = "  -11.02%  "
    "D25" = "'13.60"
    "E25" = "  -11.79%  "
    "D26" = "'0.999"
    "E26" = "  -0.06%  "
    "E27" = "  -10.87%  "
    "E28" = "  -13.63%  "
    "E29" = "  -11.42%  "
    "E30" = "  -8.99%  "
    "D31" = "'2.70"
    "E31" = "  -6.91%  "
    "E32" = "  -8.41%  "
    "D33" = "'6.64"
    "E33" = "  -17.38%  "
    "E34" = "  -13.20%  "
    "D35" = "'534.41"
    "E35" = "  -9.94%  "
    "E36" = "  +0.03%  "
    "D37" = "'0.0457"
    "E37" = "  -7.73%  "
    "D38" = "'53.41"
    "E38" = "  -5.55%  "
    "D39" = "'0.0861"
    "E39" = "  -12.26%  "
    "D40" = "'9.07"
    "E40" = "  -15.94%  "
    "E41" = "  -11.20%  "
    "D42" = "'2.81"
    "E42" = "  -14.81%  "
    "D43" = "2.947.68"
    "E43" = "  -10.92%  "
    "E44" = "  -11.60%  "
    "D45" = "0.0₃0596"
    "E45" = "  -16.33%  "
    "D46" = "'2.20"
    "E46" = "  -10.24%  "
    "D47" = "'26.93"
    "E47" = "  -14.22%  "
    "B48" = "USDe"
    "C48" = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
    "D48" = "'1.00"
    "E48" = "  -0.04%  "
    "B49" = "ThetaToken"
    "C49" = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
    "D49" = "'2.36"
    "E49" = "  -15.80%  "
    "D50" = "'125.13"
    "E50" = "  -6.49%  "
    "E51" = "  -11.14%  "
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

$wb.Save()